# Get_Location_FullyDressedUseCase.docx edit
#
# The diff does two things that, together, describe the Word "_GoBack"
# bookmark (auto-maintained at the site of the last edit) jumping from
# the end of the document to a point inside the "Use case:" line, where
# an "s" was deleted from "Locations" to make it "Location":
#
#   1. "Use case: Get Locations;" -> "Use case: Get Location;"
#      (the run is split so the bookmark can sit right before the ";")
#   2. The old _GoBack bookmark (previously sitting right after the
#      "Post-condition: N/A" paragraph) is removed.
#
# Net effect: same bookmark, relocated to just after "Get Location" and
# before the trailing ";", with the "s" removed from "Locations".

$d = $word.ActiveDocument

# --- 1. Drop the old _GoBack bookmark, wherever it currently lives ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. Fix the text: "Get Locations;" -> "Get Location;" -------------
$d.Content.Find.Execute("Use case: Get Locations;", $false, $false, $false, `
                         $false, $false, $true, 1, $false, `
                         "Use case: Get Location;", 2)

# --- 3. Re-create _GoBack right between "Get Location" and ";" --------
$r = $d.Content
$r.Find.Execute("Use case: Get Location", $false, $false, $false, $false, `
                 $false, $true, 1, $false, "", 0)
$r.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r)
